$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-detected as numbers by Excel;
# force them to remain Text so the stored cell type matches the original (string) cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.702.76"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.633.22"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "217.69"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "0.497"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "18.97"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.860.55"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.629.43"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "63.95"
$ws.Range("D17").Value = "26.677.66"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "210.91"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("E23").Value = "  -8.72%  "
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "146.71"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("D28").Value = "7.00"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "1.259.36"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("D38").Value = "0.519"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "0.795"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").Value = "0.797"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("D43").Value = "1.771.53"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").Value = "90.99"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "59.58"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.407"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -2.86%  "
